$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift each data row left by one column (drop oldest period, make room for newest)
# then fill in the new rightmost (column M) value for the newest period.

# Row 8: period-length headers
$ws.Range("E8:M8").Copy($ws.Cells.Item(8, 4))
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish-date headers
$ws.Range("E9:M9").Copy($ws.Cells.Item(9, 4))
$ws.Cells.Item(9, 9).Value = "1402-02-27 (12)"
$ws.Cells.Item(9, 13).Value = "1402-02-27 (3)"

# Row 11: Sales (فروش)
$ws.Range("E11:M11").Copy($ws.Cells.Item(11, 4))
$ws.Cells.Item(11, 13).Value = 59954

# Row 12: Cost of goods sold (بهای تمام شده کالای فروش رفته)
$ws.Range("E12:M12").Copy($ws.Cells.Item(12, 4))
$ws.Cells.Item(12, 13).Value = -25118

# Row 13: Gross profit (سود (زیان) ناخالص)
$ws.Range("E13:M13").Copy($ws.Cells.Item(13, 4))
$ws.Cells.Item(13, 13).Value = 34837

# Row 14: General & admin expenses (هزینه های عمومی)
$ws.Range("E14:M14").Copy($ws.Cells.Item(14, 4))
$ws.Cells.Item(14, 13).Value = -5120

# Row 16: Other operating income/expense
$ws.Range("E16:M16").Copy($ws.Cells.Item(16, 4))
$ws.Cells.Item(16, 13).Value = 1076

# Row 17: Operating profit
$ws.Range("E17:M17").Copy($ws.Cells.Item(17, 4))
$ws.Cells.Item(17, 13).Value = 30793

# Row 18: Finance costs
$ws.Range("E18:M18").Copy($ws.Cells.Item(18, 4))
$ws.Cells.Item(18, 13).Value = -157

# Row 19: Other non-operating income/expense
$ws.Range("E19:M19").Copy($ws.Cells.Item(19, 4))
$ws.Cells.Item(19, 13).Value = 3543

# Row 20: Profit before tax
$ws.Range("E20:M20").Copy($ws.Cells.Item(20, 4))
$ws.Cells.Item(20, 13).Value = 34178

# Row 21: Tax
$ws.Range("E21:M21").Copy($ws.Cells.Item(21, 4))
$ws.Cells.Item(21, 13).Value = -5006

# Row 22: Net profit from continuing operations
$ws.Range("E22:M22").Copy($ws.Cells.Item(22, 4))
$ws.Cells.Item(22, 13).Value = 29172

# Row 24: Net profit
$ws.Range("E24:M24").Copy($ws.Cells.Item(24, 4))
$ws.Cells.Item(24, 13).Value = 29172

# Row 26: Capital (سرمایه)
$ws.Range("E26:M26").Copy($ws.Cells.Item(26, 4))
$ws.Cells.Item(26, 13).Value = 1842
